$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Set the Jurisdiction value cell to "FRANCE"
$ws.Range("B11").Value = "FRANCE"

# Update the Date value cell with the new timestamp
$ws.Range("B8").Value = "2025-07-11T12:29:53+00:00"
